# Update _CategoryLookup and _SubcategoryLookup sheets so the single
# lookup column shows 'Display Name (code)' instead of just the raw code.
# Row 1 holds the header (renamed to *_display); rows 2+ hold the values.

$wb = $excel.ActiveWorkbook

$categoryValues = @(
    "category_display",
    "Income (income)",
    "Transfers In (transfers_in)",
    "Transfers Out (transfers_out)",
    "Loan Payments (EMIs) (loans_payments)",
    "Regular Investments / Commitments (investments_commitments)",
    "Insurance Premiums (insurance_premiums)",
    "Housing (Rent & Society) (housing_fixed)",
    "Utilities (utilities)",
    "Entertainment & OTT (entertainment)",
    "Food & Dining / Nightlife (food_dining)",
    "Groceries (groceries)",
    "Medical & Healthcare (medical)",
    "Fitness & Sports (fitness)",
    "Transport & Travel (transport)",
    "Shopping & Retail (shopping)",
    "Education (education)",
    "Child Care (child_care)",
    "Motor Maintenance (motor_maintenance)",
    "Pets (pets)",
    "Bank Interest & Fees (banks)",
    "Government Taxes (govt_tax)",
    "Charity & Donations (charity_donations)",
    "Festivals, Rituals & Celebrations (festivals_rituals)",
    "Family Support & Obligations (family_support)",
    "Business & Freelance Expenses (business_expenses)",
    "Government Benefits & Subsidies (govt_benefits)",
    "Credit Cards (credit_cards)"
)

$subcategoryValues = @(
    "subcategory_display",
    "Online Shopping (amazon)",
    "Bank Charges / Fees (bank_charges)",
    "Interest Credit (Bank) (bank_interest)",
    "Savings Sweep / Auto (bank_sweep)",
    "Ads, Marketing, Campaigns (biz_marketing)",
    "Other Business Expenses (biz_other)",
    "Materials / Inventory (biz_raw_materials)",
    "Office / Shop Rent (biz_rent)",
    "Salaries/Stipends to Staff (biz_salary)",
    "Tools, SaaS, Domains, Hosting (biz_tools_software)",
    "Client Travel & Stays (biz_travel)",
    "Credit Card Balance Transfer (cc_balance_transfer)",
    "Credit Card Bill Payment (cc_bill_payment)",
    "Credit Card Cashback / Rewards (Income) (cc_cashback)",
    "Credit Card EMI / Installments (cc_emi)",
    "Credit Card Fees (Annual/Late/Overlimit) (cc_fees)",
    "Foreign Transaction Fee (cc_foreign_transaction)",
    "Credit Card Interest / Finance Charges (cc_interest)",
    "Other Credit Card Charges (cc_other)",
    "Church Offerings (char_church)",
    "Crowdfunding / Medical Help (char_crowdfund)",
    "Gurudwara / Langar / Seva (char_gurudwara)",
    "Mosque / Zakat / Sadaqah (char_mosque)",
    "NGOs & Social Causes (char_ngo)",
    "Other Donations & Offerings (char_other)",
    "Temple / Puja / Hundi (char_temple)",
    "Daycare / Babysitting (child_daycare)",
    "Child Education Expenses (child_education)",
    "Online Courses / Certifications (edu_online)",
    "School / College Fees (edu_school_fees)",
    "Tuition / Coaching (edu_tuition)",
    "Electricity Bill (electricity)",
    "Adventure & Recreation (ent_adventure)",
    "Amusement Parks & Events (ent_amusement)",
    "Casinos & Gambling (ent_casinos)",
    "Cultural & Festive Events (ent_cultural)",
    "Video Games & eSports (ent_gaming)",
    "Movies & TV / OTT (ent_movies_ott)",
    "Museums & Art Exhibitions (ent_museums_arts)",
    "Music & Audio (ent_music)",
    "Nightlife & Parties (ent_nightlife)",
    "Other Entertainment (ent_other)",
    "Sporting Events & Tickets (ent_sports_events)",
    "Paying Relatives' Education Fees (fam_education_help)",
    "Helping Family with Medical Bills (fam_medical_help)",
    "Other Family Support (fam_other)",
    "Monthly Support to Parents/In-laws (fam_parents)",
    "Support to Siblings/Relatives (fam_relatives)",
    "Cafés & Bakeries (fd_cafes)",
    "Desserts & Sweet Shops (fd_desserts)",
    "Fine & Casual Dining (fd_fine)",
    "Online Food Delivery (Swiggy/Zomato) (fd_online)",
    "Pan / Cigarette Shop (fd_pan_shop)",
    "Pubs & Bars (fd_pubs_bars)",
    "Quick Service / Fast Food (fd_quick_service)",
    "Street Food & Local Eateries (fd_street_food)",
    "Festival Clothing & Wear (fest_clothing)",
    "Decorations & Lights (fest_decor)",
    "Festival Gifts (fest_gifts)",
    "Other Festival Expenses (fest_other)",
    "Pandal / Mandap & Events (fest_pandal)",
    "Puja Items & Ritual Needs (fest_pooja_items)",
    "Sweets & Snacks (fest_sweets)",
    "Gyms & Fitness Centers (fit_gyms)",
    "Sports & Coaching / Gear (fit_sports)",
    "Flights (flight)",
    "Other Government Benefits (govt_other_benefits)",
    "PM Kisan, Pensions, Welfare Schemes (govt_scheme)",
    "LPG, Fertilizer & Other Subsidies (govt_subsidy)",
    "Vegetable & Fruit Stores (groc_fv)",
    "Hypermarkets / Department Stores (groc_hyper)",
    "Meat / Poultry / Seafood (groc_meat)",
    "Online Groceries / Q-commerce (groc_online)",
    "Maid / Security / House Help (house_maid)",
    "Rent (house_rent)",
    "Society / Maintenance (house_society)",
    "Business Income (inc_business)",
    "Dividends (inc_dividend)",
    "Interest Income (inc_interest)",
    "Other Income (inc_other)",
    "Salary / Payroll (inc_salary)",
    "Side Income / Freelance (inc_side)",
    "Tax Refund (inc_tax_refund)",
    "Health Insurance (ins_health)",
    "Home / Other Insurance (ins_home_other)",
    "Life Insurance (ins_life)",
    "Motor Insurance (ins_motor)",
    "FD / RD (inv_fd_rd)",
    "Gold / SGB (inv_gold)",
    "NPS (inv_nps)",
    "PPF (inv_ppf)",
    "Mutual Fund SIP (inv_sip)",
    "Stocks / ETFs (inv_stocks)",
    "Two-Wheeler EMI (loan_bike)",
    "Car Loan EMI (loan_car)",
    "Credit Card Bill Payment (loan_cc_bill)",
    "Home Loan EMI (loan_home)",
    "Other Loan/Debt Payment (loan_other)",
    "Personal Loan EMI (loan_personal)",
    "Education Loan EMI (loan_student)",
    "Medical Apps / Services (med_apps)",
    "Dental Care (med_dental)",
    "Eye Care / Optometry (med_eye)",
    "General / Hospitals / Nursing (med_general)",
    "Other Medical (med_other)",
    "Pharmacies & Supplements (med_pharma)",
    "Mobile/Internet (mobile)",
    "Motor Insurance (motor_insurance)",
    "General Services / Repairs (motor_services)",
    "Pet Food (pet_food)",
    "Grooming / Boarding / Bathing (pet_grooming)",
    "Pet Insurance (pet_insurance)",
    "Vaccination / Vet (pet_vaccine)",
    "Recurring Deposit (recurring_deposit)",
    "Restaurant (restaurant)",
    "Automotive Supplies (shop_auto_supplies)",
    "Beauty & Personal Care (shop_beauty)",
    "Books & Media (shop_books_media)",
    "Children & Toys (shop_children_toys)",
    "Clothing & Accessories (shop_clothing)",
    "Electronics & Gadgets (shop_electronics)",
    "General Merchandise (shop_general)",
    "Gifts & Novelties (shop_gifts)",
    "Hobbies & Crafts (shop_hobbies)",
    "Home & Kitchen / Furnishings (shop_home_kitchen)",
    "Online Marketplaces (shop_marketplaces)",
    "Pet Supplies (shop_pet_supplies)",
    "Sports & Outdoor (shop_sports_outdoor)",
    "Stationery & Office Supplies (shop_stationery)",
    "Supermarket (supermarket)",
    "GST / Challan (tax_gst)",
    "Income Tax / TDS (tax_income)",
    "Other Government Taxes (tax_other)",
    "Transport Apps (Uber/Ola etc.) (tr_apps)",
    "Fuel / Petrol / Diesel (tr_fuel)",
    "Cash/Cheque/ATM Deposit (tr_in_deposit)",
    "Account Transfer In (tr_in_internal)",
    "From Investments / PF / Retirement (tr_in_invest_ret)",
    "Other Transfer In (tr_in_other)",
    "Savings Sweep / Sweep In (tr_in_savings)",
    "Hotels / Stays / Airbnb (tr_lodging)",
    "Other Transport (tr_other)",
    "Cash / ATM Withdrawal (tr_out_atm)",
    "Other Transfer Out (tr_out_other)",
    "Transfer to Savings/Joint (tr_out_savings)",
    "Sweep Out (FD / Linked) (tr_out_sweep)",
    "Bank → Wallet / UPI (tr_out_wallet)",
    "Public Transit (Rail/Metro/Bus) (tr_public)",
    "Taxis / Auto / Ride-share (tr_taxis)",
    "Tolls / FASTag (tr_tolls)",
    "Flights / Bus / Train / Cruise (tr_travel)",
    "Cabs (uber)",
    "Internet / Broadband (util_broadband)",
    "DTH / Cable TV (util_dth_cable)",
    "Electricity (util_electricity)",
    "Gas / LPG / PNG (util_gas_lpg)",
    "Mobile / Telephone (util_mobile)",
    "Other Utilities (util_other)",
    "Sewage / Waste Mgmt (util_sewage_waste)",
    "Water (util_water)",
    "Water Bill (water)",
    "Food Delivery (zomato)"
)

$wsCategory = $wb.Worksheets.Item("_CategoryLookup")
for ($i = 0; $i -lt $categoryValues.Length; $i++) {
    $wsCategory.Cells.Item($i + 1, 1).Value = $categoryValues[$i]
}

$wsSubcategory = $wb.Worksheets.Item("_SubcategoryLookup")
for ($i = 0; $i -lt $subcategoryValues.Length; $i++) {
    $wsSubcategory.Cells.Item($i + 1, 1).Value = $subcategoryValues[$i]
}

